$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''67.552.14'
$ws.Range("E2").Value = '  -1.08%  '
$ws.Range("D3").Value = '''3.773.09'
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").Value = '''595.00'
$ws.Range("E5").Value = '  -0.02%  '
$ws.Range("D6").Value = '''166.21'
$ws.Range("E6").Value = '  -0.57%  '
$ws.Range("D7").Value = '''3.772.03'
$ws.Range("E7").Value = '  +0.72%  '
$ws.Range("E8").Value = '  +0.06%  '
$ws.Range("E9").Value = '  -0.15%  '
$ws.Range("D10").Value = '''0.159'
$ws.Range("E10").Value = '  -0.05%  '
$ws.Range("E11").Value = '  -1.93%  '
$ws.Range("D12").Value = '''0.448'
$ws.Range("E12").Value = '  +0.26%  '
$ws.Range("E13").Value = '  -1.61%  '
$ws.Range("D14").Value = '''36.19'
$ws.Range("E14").Value = '  +0.52%  '
$ws.Range("D15").Value = '''4.404.22'
$ws.Range("E15").Value = '  +0.67%  '
$ws.Range("D16").Value = '''3.789.76'
$ws.Range("E16").Value = '  +0.98%  '
$ws.Range("D17").Value = '''18.42'
$ws.Range("E17").Value = '  +2.90%  '
$ws.Range("D18").Value = '''67.531.17'
$ws.Range("E19").Value = '  +0.17%  '
$ws.Range("E20").Value = '  -0.16%  '
$ws.Range("D21").Value = '''10.05'
$ws.Range("E21").Value = '  -6.58%  '
$ws.Range("D22").Value = '''456.12'
$ws.Range("E22").Value = '  -2.01%  '
$ws.Range("E24").Value = '  +7.39%  '
$ws.Range("D25").Value = '''83.28'
$ws.Range("E25").Value = '  -1.36%  '
$ws.Range("D26").Value = '''11.90'
$ws.Range("E26").Value = '  -0.74%  '
$ws.Range("D28").Value = '''10.07'
$ws.Range("E28").Value = '  +0.05%  '
$ws.Range("E29").Value = '  +0.06%  '
$ws.Range("E30").Value = '  -0.02%  '
$ws.Range("D31").Value = '''7.26'
$ws.Range("E31").Value = '  -0.46%  '
$ws.Range("D32").Value = '''29.73'
$ws.Range("E32").Value = '  -0.33%  '
$ws.Range("D33").Value = '''2.18'
$ws.Range("E33").Value = '  +0.51%  '
$ws.Range("E34").Value = '  +0.01%  '
$ws.Range("D35").Value = '''1.00'
$ws.Range("E35").Value = '  -0.08%  '
$ws.Range("D36").Value = '''3.723.23'
$ws.Range("E36").Value = '  +0.57%  '
$ws.Range("E37").Value = '  -0.79%  '
$ws.Range("E38").Value = '  -0.88%  '
$ws.Range("E39").Value = '  -1.03%  '
$ws.Range("D40").Value = '''0.996'
$ws.Range("E40").Value = '  -0.21%  '
$ws.Range("E41").Value = '  -0.91%  '
$ws.Range("D42").Value = '''0.999'
$ws.Range("E42").Value = '  -0.04%  '
$ws.Range("E43").Value = '  -0.01%  '
$ws.Range("D44").Value = '''45.05'
$ws.Range("E44").Value = '  +5.11%  '
$ws.Range("E45").Value = '  -1.56%  '
$ws.Range("D46").Value = '''47.15'
$ws.Range("E46").Value = '  +2.84%  '
$ws.Range("E47").Value = '  -2.64%  '
$ws.Range("D48").Value = '''148.04'
$ws.Range("E48").Value = '  +0.93%  '
$ws.Range("E49").Value = '  -4.17%  '
$ws.Range("D50").Value = '''390.30'
$ws.Range("E50").Value = '  +0.43%  '
$ws.Range("D51").Value = '''25.55'
$ws.Range("E51").Value = '  -0.32%  '
